# The document's first-page header/footer carry a BTEC logo (header) and two
# Pearson logos (primary + first-page footers). This edit simply renames the
# embedded picture objects (the wp:docPr / pic:cNvPr "name" attribute used
# internally by Word for each inline picture), swapping:
#   footer PearsonLogo pictures: image1.png -> image2.png
#   header BTec_Logo-Orange picture: image2.jpg -> image1.jpg
# No visual/content change - this is a rename of the picture object names.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Primary (default) footer: PearsonLogo, currently "image1.png" ---
$ftrPrimary = $sec.Footers.Item(1)
if ($ftrPrimary.Range.InlineShapes.Count -ge 1) {
    $ftrPrimary.Range.InlineShapes.Item(1).Name = "image2.png"
}

# --- First-page footer: PearsonLogo, currently "image1.png" ---
$ftrFirst = $sec.Footers.Item(2)
if ($ftrFirst.Range.InlineShapes.Count -ge 1) {
    $ftrFirst.Range.InlineShapes.Item(1).Name = "image2.png"
}

# --- First-page header: BTec_Logo-Orange, currently "image2.jpg" ---
$hdrFirst = $sec.Headers.Item(2)
if ($hdrFirst.Range.InlineShapes.Count -ge 1) {
    $hdrFirst.Range.InlineShapes.Item(1).Name = "image1.jpg"
}

Write-Output "Renamed header/footer logo picture objects."
